$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the existing D:K data
# (and the blank-but-styled cells) one column to the right, into E:L, and
# creates a brand-new (blank, default-styled) column D.
$ws.Columns("D:D").Insert()

# The freshly inserted column D doesn't inherit the number formatting that
# columns D:K previously shared. Copy that formatting from the (now shifted)
# column E so the new column D matches (date format on the header rows,
# number format everywhere else).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D with the new period's figures for every row
# that carries data in this table (three stacked statements: Income
# Statement, Balance Sheet, Cash Flow Statement).

# -- Income Statement --
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 406500
$ws.Range("D9").Value = 275100
$ws.Range("D10").Value = 131400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 89300
$ws.Range("D15").Value = 58300
$ws.Range("D17").Value = 522100
$ws.Range("D18").Value = -115600
$ws.Range("D20").Value = -7500
$ws.Range("D21").Value = -71700
$ws.Range("D22").Value = 30400
$ws.Range("D23").Value = -153500
$ws.Range("D24").Value = 18300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -171800
$ws.Range("D27").Value = -171500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 7500
$ws.Range("D33").Value = -171500
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -171500

# -- Balance Sheet --
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 371800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 244200
$ws.Range("D44").Value = 29500
$ws.Range("D45").Value = 37800
$ws.Range("D46").Value = 683300
$ws.Range("D47").Value = 1000
$ws.Range("D48").Value = 1089900
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 53500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1827700
$ws.Range("D57").Value = 31900
$ws.Range("D58").Value = 8600
$ws.Range("D59").Value = 117800
$ws.Range("D60").Value = 158400
$ws.Range("D61").Value = 430400
$ws.Range("D62").Value = 94000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 683900
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -210800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1143800
$ws.Range("D77").Value = 0

# -- Cash Flow Statement --
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -171500
$ws.Range("D83").Value = 51300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 3900
$ws.Range("D91").Value = -21400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 68500
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -128100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -55600
